$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Approved" for rows 3-5 in column I (the "Approved/Rejected" column),
# matching the value already present in I2.
$ws.Range("I3").Value = "Approved"
$ws.Range("I4").Value = "Approved"
$ws.Range("I5").Value = "Approved"

# Reflect the scrolled viewport / new active selection recorded in the sheet
# (user scrolled right to column E and selected G14).
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G14").Select()
